# Weekly update: insert two new rows (newer market data) right before the
# last existing block of "Packham's Triumph" entries for
# "Terminal Hortofrutícola Agro Chillán - Pera", pushing the rest of the
# rows down by two (old row 604 -> new row 606, ..., old row 639 -> new row 641).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 604-605; everything from the old row 604 onward
# shifts down to make room.
$ws.Rows("604:605").Insert()

# New row 604: Packham's Triumph / Primera, newest week (2023-12-07 = serial 45267)
$ws.Cells.Item(604, 1).Value = 7
$ws.Cells.Item(604, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(604, 3).Value = "Ñuble"
$ws.Cells.Item(604, 4).Value = 45267
$ws.Cells.Item(604, 5).Value = 16
$ws.Cells.Item(604, 6).Value = "Fruta"
$ws.Cells.Item(604, 7).Value = 100104
$ws.Cells.Item(604, 8).Value = "Frutos de pepita"
$ws.Cells.Item(604, 9).Value = 100104005
$ws.Cells.Item(604, 10).Value = "Pera"
$ws.Cells.Item(604, 11).Value = "Packham's Triumph"
$ws.Cells.Item(604, 12).Value = "Primera"
$ws.Cells.Item(604, 13).Value = 180
$ws.Cells.Item(604, 14).Value = 15000
$ws.Cells.Item(604, 15).Value = 15000
$ws.Cells.Item(604, 16).Value = 15000
$ws.Cells.Item(604, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(604, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(604, 19).Value = 833
$ws.Cells.Item(604, 20).Value = 18

# New row 605: Packham's Triumph / Segunda, same newest week
$ws.Cells.Item(605, 1).Value = 7
$ws.Cells.Item(605, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(605, 3).Value = "Ñuble"
$ws.Cells.Item(605, 4).Value = 45267
$ws.Cells.Item(605, 5).Value = 16
$ws.Cells.Item(605, 6).Value = "Fruta"
$ws.Cells.Item(605, 7).Value = 100104
$ws.Cells.Item(605, 8).Value = "Frutos de pepita"
$ws.Cells.Item(605, 9).Value = 100104005
$ws.Cells.Item(605, 10).Value = "Pera"
$ws.Cells.Item(605, 11).Value = "Packham's Triumph"
$ws.Cells.Item(605, 12).Value = "Segunda"
$ws.Cells.Item(605, 13).Value = 100
$ws.Cells.Item(605, 14).Value = 13000
$ws.Cells.Item(605, 15).Value = 13000
$ws.Cells.Item(605, 16).Value = 13000
$ws.Cells.Item(605, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(605, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(605, 19).Value = 722
$ws.Cells.Item(605, 20).Value = 18
